$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reshape the data rows -------------------------------------------------
# Original sheet has 4 data rows (Excel rows 2-5, ids 163,164,165,166).
# Target sheet keeps only 2 data rows:
#   row 2 -> a brand new record (id 167, subject "英语", title "yingyu", ...)
#   row 3 -> the original row 2 record (id 163) moved down, untouched
#
# Remove rows 3..5 (ids 164,165,166) first, then insert a fresh blank row
# above the remaining data row so the original id=163 record slides from
# row 2 down to row 3.
$ws.Rows("3:5").Delete()
$ws.Rows("2:2").Insert()

# Copy the formatting (styles, number formats) from row 3 (the preserved
# id=163 record) onto the newly inserted blank row 2 so it matches the
# rest of the table.
$ws.Range("A3:T3").Copy()
$ws.Range("A2:T2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows("2:2").RowHeight = 180

# --- Populate the new row 2 record (id 167) --------------------------------
$ws.Range("A2").Value2 = 167
$ws.Range("B2").Value2 = 1
$ws.Range("C2").Value2 = "英语"
$ws.Range("D2").Value2 = "yingyu"
$ws.Range("E2").Value2 = ""
$ws.Range("F2").Value2 = ""
$ws.Range("G2").Value2 = ""
$ws.Range("H2").Value2 = $true
$ws.Range("I2").Value2 = 1
$ws.Range("J2").Value2 = 1
$ws.Range("K2").Value2 = "未分类"
$ws.Range("L2").Value2 = $true
$ws.Range("M2").Value2 = ""
$ws.Range("N2").Value2 = 2
$ws.Range("O2").Value2 = "d41d8cd98f00b204e9800998ecf8427e"
$ws.Range("P2").Value2 = ""
$ws.Range("Q2").Value2 = $true
$ws.Range("R2").Value2 = "ben"
$ws.Range("S2").Value2 = 45223
$ws.Range("T2").Value2 = 45223

$ws.Range("A1").Select()

# --- Fix up the pictures anchored next to the data rows --------------------
# Before: 3 pictures anchored at rows 2, 4 and 5 (Image 1/2/3).
# After: only the first picture remains, now anchored next to row 3
# (the relocated id=163 record), i.e. moved down by one row (180pt).
$ws.Shapes.Item("Image 3").Delete()
$ws.Shapes.Item("Image 2").Delete()
$ws.Shapes.Item("Image 1").Top = $ws.Rows("3:3").Top
